# Apply updated production figures for the new RAAL model.
$wb = $excel.ActiveWorkbook

# --- Daily sheet ---
$daily = $wb.Worksheets.Item("Daily")
$daily.Range("G2").Value = 2896.43
$daily.Range("H2").Value = 6152.96
$daily.Range("I2").Value = 718.61
$daily.Range("J2").Value = 724.11
$daily.Range("L2").Value = 724.11

# --- Hourly sheet ---
$hourly = $wb.Worksheets.Item("Hourly")

$hourly.Range("I9").Value = 26.35

$hourly.Range("H10").Value = 90.56999999999999
$hourly.Range("I10").Value = 392.99
$hourly.Range("J10").Value = 44.69
$hourly.Range("K10").Value = 22.64
$hourly.Range("M10").Value = 22.64

$hourly.Range("H11").Value = 233.23
$hourly.Range("I11").Value = 620.88
$hourly.Range("J11").Value = 70.06

$hourly.Range("H12").Value = 357.93
$hourly.Range("I12").Value = 726.28
$hourly.Range("J12").Value = 84.45999999999999
$hourly.Range("K12").Value = 89.48
$hourly.Range("M12").Value = 89.48

$hourly.Range("H13").Value = 443.53
$hourly.Range("I13").Value = 778.59
$hourly.Range("J13").Value = 92.47
$hourly.Range("K13").Value = 110.88
$hourly.Range("M13").Value = 110.88

$hourly.Range("H14").Value = 479.44
$hourly.Range("I14").Value = 797.47
$hourly.Range("J14").Value = 95.5
$hourly.Range("K14").Value = 119.86
$hourly.Range("M14").Value = 119.86

$hourly.Range("H15").Value = 461.59
$hourly.Range("I15").Value = 788.38
$hourly.Range("J15").Value = 93.98999999999999

$hourly.Range("H16").Value = 391.95
$hourly.Range("I16").Value = 748.65
$hourly.Range("J16").Value = 87.75

$hourly.Range("H17").Value = 278.8
$hourly.Range("I17").Value = 665
$hourly.Range("J17").Value = 75.79000000000001
$hourly.Range("K17").Value = 69.7
$hourly.Range("M17").Value = 69.7

$hourly.Range("H18").Value = 138.64
$hourly.Range("I18").Value = 493.97
$hourly.Range("J18").Value = 55.13

$hourly.Range("I19").Value = 114.4
